$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "정률할인"
$ws.Range("D3").Value = "정률할인"
$ws.Range("D4").Value = "정액할인"
$ws.Range("D5").Value = "정액할인"
$ws.Range("D6").Value = "수량별 정액할인"
$ws.Range("D7").Value = "정액할인"
$ws.Range("D8").Value = "정률할인"
